$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value = 1000.75  # H31: 0 -> 1000.75
$ws.Cells.Item(31, 9).Value = 1001  # I31: 0 -> 1001
$ws.Cells.Item(31, 10).Value = 1000  # J31: 0 -> 1000
$ws.Cells.Item(31, 11).Value = 3003  # K31: 0 -> 3003
$ws.Cells.Item(31, 12).Value = 3000  # L31: 0 -> 3000
$ws.Cells.Item(31, 13).Value = -2773  # M31: None -> -2773
$ws.Cells.Item(31, 14).Value = -3460  # N31: None -> -3460

$ws.Cells.Item(43, 8).Value = 1329.4286  # H43: 1620.25 -> 1329.4286
$ws.Cells.Item(43, 9).Value = 0  # I43: 2000 -> 0
$ws.Cells.Item(43, 10).Value = 1329.4286  # J43: 1566 -> 1329.4286
$ws.Cells.Item(43, 11).Value = 0  # K43: 2000 -> 0
$ws.Cells.Item(43, 12).Value = 1329.4286  # L43: 1566 -> 1329.4286
$ws.Cells.Item(43, 13).ClearContents()  # M43: -1931 -> (removed)
$ws.Cells.Item(43, 14).Value = -1467.4286  # N43: -1704 -> -1467.4286

$ws.Cells.Item(116, 8).Value = 5526.615  # H116: 4455.75 -> 5526.615
$ws.Cells.Item(116, 9).Value = 3420.8572  # I116: 2000 -> 3420.8572
$ws.Cells.Item(116, 10).Value = 7983.3335  # J116: 4619.467 -> 7983.3335
$ws.Cells.Item(116, 11).Value = 3420.8572  # K116: 2000 -> 3420.8572
$ws.Cells.Item(116, 12).Value = 7983.3335  # L116: 4619.467 -> 7983.3335
$ws.Cells.Item(116, 13).Value = 21.14280000000008  # M116: 1442 -> 21.14280000000008
$ws.Cells.Item(116, 14).Value = -14867.3335  # N116: -11503.467 -> -14867.3335

$ws.Cells.Item(127, 8).Value = 847.875  # H127: 1325 -> 847.875
$ws.Cells.Item(127, 9).Value = 591.5  # I127: 659.8 -> 591.5
$ws.Cells.Item(127, 10).Value = 1617  # J127: 1694.5555 -> 1617
$ws.Cells.Item(127, 11).Value = 1774.5  # K127: 1979.4 -> 1774.5
$ws.Cells.Item(127, 12).Value = 4851  # L127: 5083.666499999999 -> 4851
$ws.Cells.Item(127, 13).Value = 3185.5  # M127: 2980.6 -> 3185.5
$ws.Cells.Item(127, 14).Value = -14771  # N127: -15003.6665 -> -14771

$ws.Cells.Item(128, 8).Value = 500022300  # H128: 500024000 -> 500022300
$ws.Cells.Item(128, 10).Value = 500022300  # J128: 500024000 -> 500022300
$ws.Cells.Item(128, 12).Value = 500022300  # L128: 500024000 -> 500022300
$ws.Cells.Item(128, 14).Value = -500032260  # N128: -500033960 -> -500032260

$ws.Cells.Item(132, 8).Value = 6374.846  # H132: 3992.9404 -> 6374.846
$ws.Cells.Item(132, 9).Value = 7349.76  # I132: 3155.9836 -> 7349.76
$ws.Cells.Item(132, 10).Value = 5472.148  # J132: 6212.696 -> 5472.148
$ws.Cells.Item(132, 11).Value = 22049.28  # K132: 9467.950800000001 -> 22049.28
$ws.Cells.Item(132, 12).Value = 16416.444  # L132: 18638.088 -> 16416.444
$ws.Cells.Item(132, 13).Value = -19519.28  # M132: -6937.950800000001 -> -19519.28
$ws.Cells.Item(132, 14).Value = -21476.444  # N132: -23698.088 -> -21476.444

$ws.Cells.Item(138, 8).Value = 1428.6522  # H138: 1269.5416 -> 1428.6522
$ws.Cells.Item(138, 9).Value = 1105.317  # I138: 1015.9778 -> 1105.317
$ws.Cells.Item(138, 10).Value = 4080  # J138: 5073 -> 4080
$ws.Cells.Item(138, 11).Value = 3315.951  # K138: 3047.9334 -> 3315.951
$ws.Cells.Item(138, 12).Value = 12240  # L138: 15219 -> 12240
$ws.Cells.Item(138, 13).Value = 1824.049  # M138: 2092.0666 -> 1824.049
$ws.Cells.Item(138, 14).Value = -22520  # N138: -25499 -> -22520

$ws.Cells.Item(141, 8).Value = 4904.372  # H141: 4271.574 -> 4904.372
$ws.Cells.Item(141, 9).Value = 1192.7354  # I141: 1084.7632 -> 1192.7354
$ws.Cells.Item(141, 10).Value = 18926.111  # J141: 11840.25 -> 18926.111
$ws.Cells.Item(141, 11).Value = 3578.2062  # K141: 3254.2896 -> 3578.2062
$ws.Cells.Item(141, 12).Value = 56778.333  # L141: 35520.75 -> 56778.333
$ws.Cells.Item(141, 13).Value = 1601.7938  # M141: 1925.7104 -> 1601.7938
$ws.Cells.Item(141, 14).Value = -67138.333  # N141: -45880.75 -> -67138.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9203.489  # H32: 10123.275 -> 9203.489
$ws.Cells.Item(32, 9).Value = 4622.1357  # I32: 5160.875 -> 4622.1357
$ws.Cells.Item(32, 11).Value = 4622.1357  # K32: 5160.875 -> 4622.1357
$ws.Cells.Item(32, 13).Value = -4335.1357  # M32: -4873.875 -> -4335.1357

$ws.Cells.Item(61, 8).Value = 1131.8833  # H61: 1443.5927 -> 1131.8833
$ws.Cells.Item(61, 9).Value = 857.5526  # I61: 1302.1818 -> 857.5526
$ws.Cells.Item(61, 10).Value = 1605.7273  # J61: 1665.8096 -> 1605.7273
$ws.Cells.Item(61, 11).Value = 857.5526  # K61: 1302.1818 -> 857.5526
$ws.Cells.Item(61, 12).Value = 1605.7273  # L61: 1665.8096 -> 1605.7273
$ws.Cells.Item(61, 13).Value = -645.5526  # M61: -1090.1818 -> -645.5526
$ws.Cells.Item(61, 14).Value = -2029.7273  # N61: -2089.8096 -> -2029.7273

$ws.Cells.Item(132, 8).Value = 2020397.2  # H132: 2722942 -> 2020397.2
$ws.Cells.Item(132, 9).Value = 3276.9285  # I132: 4900.3887 -> 3276.9285
$ws.Cells.Item(132, 10).Value = 3681555  # J132: 4470254.5 -> 3681555
$ws.Cells.Item(132, 11).Value = 9830.7855  # K132: 14701.1661 -> 9830.7855
$ws.Cells.Item(132, 12).Value = 11044665  # L132: 13410763.5 -> 11044665
$ws.Cells.Item(132, 13).Value = -7300.7855  # M132: -12171.1661 -> -7300.7855
$ws.Cells.Item(132, 14).Value = -11049725  # N132: -13415823.5 -> -11049725

$ws.Cells.Item(136, 8).Value = 1131.8833  # H136: 1443.5927 -> 1131.8833
$ws.Cells.Item(136, 9).Value = 857.5526  # I136: 1302.1818 -> 857.5526
$ws.Cells.Item(136, 10).Value = 1605.7273  # J136: 1665.8096 -> 1605.7273
$ws.Cells.Item(136, 11).Value = 2572.6578  # K136: 3906.5454 -> 2572.6578
$ws.Cells.Item(136, 12).Value = 4817.1819  # L136: 4997.4288 -> 4817.1819
$ws.Cells.Item(136, 13).Value = -22.65779999999995  # M136: -1356.5454 -> -22.65779999999995
$ws.Cells.Item(136, 14).Value = -9917.1819  # N136: -10097.4288 -> -9917.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1694.5  # H22: 1577.3414 -> 1694.5
$ws.Cells.Item(22, 9).Value = 1293.4  # I22: 1430.375 -> 1293.4
$ws.Cells.Item(22, 10).Value = 3700  # J22: 1784.8235 -> 3700
$ws.Cells.Item(22, 11).Value = 1293.4  # K22: 1430.375 -> 1293.4
$ws.Cells.Item(22, 12).Value = 3700  # L22: 1784.8235 -> 3700
$ws.Cells.Item(22, 13).Value = -1120.4  # M22: -1257.375 -> -1120.4
$ws.Cells.Item(22, 14).Value = -4046  # N22: -2130.8235 -> -4046

$ws.Cells.Item(134, 8).Value = 4891.8433  # H134: 5674.3955 -> 4891.8433
$ws.Cells.Item(134, 9).Value = 2195.0454  # I134: 2905.5334 -> 2195.0454
$ws.Cells.Item(134, 10).Value = 6937.6895  # J134: 7157.7144 -> 6937.6895
$ws.Cells.Item(134, 11).Value = 6585.1362  # K134: 8716.600199999999 -> 6585.1362
$ws.Cells.Item(134, 12).Value = 20813.0685  # L134: 21473.1432 -> 20813.0685
$ws.Cells.Item(134, 13).Value = -4050.1362  # M134: -6181.600199999999 -> -4050.1362
$ws.Cells.Item(134, 14).Value = -25883.0685  # N134: -26543.1432 -> -25883.0685

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 5449  # H99: 6000 -> 5449
$ws.Cells.Item(99, 9).Value = 5648  # I99: 6000 -> 5648
$ws.Cells.Item(99, 10).Value = 5250  # J99: 0 -> 5250
$ws.Cells.Item(99, 11).Value = 5648  # K99: 6000 -> 5648
$ws.Cells.Item(99, 12).Value = 5250  # L99: 0 -> 5250
$ws.Cells.Item(99, 13).Value = -4150  # M99: -4502 -> -4150
$ws.Cells.Item(99, 14).Value = -8246  # N99: None -> -8246

$ws.Cells.Item(126, 8).Value = 5449  # H126: 6000 -> 5449
$ws.Cells.Item(126, 9).Value = 5648  # I126: 6000 -> 5648
$ws.Cells.Item(126, 10).Value = 5250  # J126: 0 -> 5250
$ws.Cells.Item(126, 11).Value = 16944  # K126: 18000 -> 16944
$ws.Cells.Item(126, 12).Value = 15750  # L126: 0 -> 15750
$ws.Cells.Item(126, 13).Value = -14474  # M126: -15530 -> -14474
$ws.Cells.Item(126, 14).Value = -20690  # N126: None -> -20690

$ws.Cells.Item(132, 8).Value = 2528.45  # H132: 2298.3865 -> 2528.45
$ws.Cells.Item(132, 9).Value = 1824.6666  # I132: 1401.96 -> 1824.6666
$ws.Cells.Item(132, 10).Value = 3104.2727  # J132: 3477.8948 -> 3104.2727
$ws.Cells.Item(132, 11).Value = 5473.9998  # K132: 4205.88 -> 5473.9998
$ws.Cells.Item(132, 12).Value = 9312.8181  # L132: 10433.6844 -> 9312.8181
$ws.Cells.Item(132, 13).Value = -2943.9998  # M132: -1675.88 -> -2943.9998
$ws.Cells.Item(132, 14).Value = -14372.8181  # N132: -15493.6844 -> -14372.8181

$ws.Cells.Item(134, 8).Value = 1459.1459  # H134: 1051.013 -> 1459.1459
$ws.Cells.Item(134, 9).Value = 853.2759  # I134: 637.6799999999999 -> 853.2759
$ws.Cells.Item(134, 10).Value = 2383.8948  # J134: 1816.4445 -> 2383.8948
$ws.Cells.Item(134, 11).Value = 2559.8277  # K134: 1913.04 -> 2559.8277
$ws.Cells.Item(134, 12).Value = 7151.6844  # L134: 5449.333500000001 -> 7151.6844
$ws.Cells.Item(134, 13).Value = -24.82769999999982  # M134: 621.96 -> -24.82769999999982
$ws.Cells.Item(134, 14).Value = -12221.6844  # N134: -10519.3335 -> -12221.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 2953  # H22: 2883.3333 -> 2953
$ws.Cells.Item(22, 9).Value = 2895.25  # I22: 5550.5 -> 2895.25
$ws.Cells.Item(22, 10).Value = 2999.2  # J22: 2121.2856 -> 2999.2
$ws.Cells.Item(22, 11).Value = 8685.75  # K22: 16651.5 -> 8685.75
$ws.Cells.Item(22, 12).Value = 8997.599999999999  # L22: 6363.8568 -> 8997.599999999999
$ws.Cells.Item(22, 13).Value = -8516.75  # M22: -16482.5 -> -8516.75
$ws.Cells.Item(22, 14).Value = -9335.599999999999  # N22: -6701.8568 -> -9335.599999999999

$ws.Cells.Item(27, 8).Value = 2953  # H27: 2883.3333 -> 2953
$ws.Cells.Item(27, 9).Value = 2895.25  # I27: 5550.5 -> 2895.25
$ws.Cells.Item(27, 10).Value = 2999.2  # J27: 2121.2856 -> 2999.2
$ws.Cells.Item(27, 11).Value = 8685.75  # K27: 16651.5 -> 8685.75
$ws.Cells.Item(27, 12).Value = 8997.599999999999  # L27: 6363.8568 -> 8997.599999999999
$ws.Cells.Item(27, 13).Value = -8583.75  # M27: -16549.5 -> -8583.75
$ws.Cells.Item(27, 14).Value = -9201.599999999999  # N27: -6567.8568 -> -9201.599999999999

$ws.Cells.Item(49, 8).Value = 5750  # H49: 3643.2856 -> 5750
$ws.Cells.Item(49, 9).Value = 1250  # I49: 2003 -> 1250
$ws.Cells.Item(49, 10).Value = 8000  # J49: 3916.6667 -> 8000
$ws.Cells.Item(49, 11).Value = 3750  # K49: 6009 -> 3750
$ws.Cells.Item(49, 12).Value = 24000  # L49: 11750.0001 -> 24000
$ws.Cells.Item(49, 13).Value = -3594  # M49: -5853 -> -3594
$ws.Cells.Item(49, 14).Value = -24312  # N49: -12062.0001 -> -24312

$ws.Cells.Item(80, 8).Value = 1607.75  # H80: 1343.75 -> 1607.75
$ws.Cells.Item(80, 10).Value = 2660  # J80: 1606.1428 -> 2660
$ws.Cells.Item(80, 12).Value = 7980  # L80: 4818.428400000001 -> 7980
$ws.Cells.Item(80, 14).Value = -9852  # N80: -6690.428400000001 -> -9852

$ws.Cells.Item(83, 8).Value = 1607.75  # H83: 1343.75 -> 1607.75
$ws.Cells.Item(83, 10).Value = 2660  # J83: 1606.1428 -> 2660
$ws.Cells.Item(83, 12).Value = 23940  # L83: 14455.2852 -> 23940
$ws.Cells.Item(83, 14).Value = -33300  # N83: -23815.2852 -> -33300

$ws.Cells.Item(137, 8).Value = 3446.4255  # H137: 3503.239 -> 3446.4255
$ws.Cells.Item(137, 10).Value = 3685.9744  # J137: 3761.0527 -> 3685.9744
$ws.Cells.Item(137, 12).Value = 11057.9232  # L137: 11283.1581 -> 11057.9232
$ws.Cells.Item(137, 14).Value = -21257.9232  # N137: -21483.1581 -> -21257.9232

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1160310.1  # H132: 564885.25 -> 1160310.1
$ws.Cells.Item(132, 9).Value = 2780021.8  # I132: 851534.2 -> 2780021.8
$ws.Cells.Item(132, 10).Value = 3373.238  # J132: 3053.28 -> 3373.238
$ws.Cells.Item(132, 11).Value = 8340065.399999999  # K132: 2554602.6 -> 8340065.399999999
$ws.Cells.Item(132, 12).Value = 10119.714  # L132: 9159.84 -> 10119.714
$ws.Cells.Item(132, 13).Value = -8337535.399999999  # M132: -2552072.6 -> -8337535.399999999
$ws.Cells.Item(132, 14).Value = -15179.714  # N132: -14219.84 -> -15179.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(128, 8).Value = 333359840  # H128: 333359900 -> 333359840
$ws.Cells.Item(128, 10).Value = 333359840  # J128: 333359900 -> 333359840
$ws.Cells.Item(128, 12).Value = 333359840  # L128: 333359900 -> 333359840
$ws.Cells.Item(128, 14).Value = -333369800  # N128: -333369860 -> -333369800

$ws.Cells.Item(132, 8).Value = 31254104  # H132: 37041730 -> 31254104
$ws.Cells.Item(132, 9).Value = 62505616  # I132: 66672604 -> 62505616
$ws.Cells.Item(132, 10).Value = 2592.6875  # J132: 3131.9167 -> 2592.6875
$ws.Cells.Item(132, 11).Value = 187516848  # K132: 200017812 -> 187516848
$ws.Cells.Item(132, 12).Value = 7778.0625  # L132: 9395.750100000001 -> 7778.0625
$ws.Cells.Item(132, 13).Value = -187514318  # M132: -200015282 -> -187514318
$ws.Cells.Item(132, 14).Value = -12838.0625  # N132: -14455.7501 -> -12838.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 77685.8  # H46: 54474.4 -> 77685.8
$ws.Cells.Item(46, 10).Value = 77685.8  # J46: 54474.4 -> 77685.8
$ws.Cells.Item(46, 12).Value = 77685.8  # L46: 54474.4 -> 77685.8
$ws.Cells.Item(46, 14).Value = -78147.8  # N46: -54936.4 -> -78147.8

$ws.Cells.Item(126, 8).Value = 2942.375  # H126: 2193.2778 -> 2942.375
$ws.Cells.Item(126, 9).Value = 3894.818  # I126: 3054.9092 -> 3894.818
$ws.Cells.Item(126, 10).Value = 847  # J126: 839.2857 -> 847
$ws.Cells.Item(126, 11).Value = 11684.454  # K126: 9164.7276 -> 11684.454
$ws.Cells.Item(126, 12).Value = 2541  # L126: 2517.8571 -> 2541
$ws.Cells.Item(126, 13).Value = -9214.454000000002  # M126: -6694.7276 -> -9214.454000000002
$ws.Cells.Item(126, 14).Value = -7481  # N126: -7457.8571 -> -7481

$ws.Cells.Item(134, 8).Value = 77685.8  # H134: 54474.4 -> 77685.8
$ws.Cells.Item(134, 10).Value = 77685.8  # J134: 54474.4 -> 77685.8
$ws.Cells.Item(134, 12).Value = 233057.4  # L134: 163423.2 -> 233057.4
$ws.Cells.Item(134, 14).Value = -238127.4  # N134: -168493.2 -> -238127.4
